$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 4, shifting "スイーパー" (previously row 4) and
# everything below it down by two rows.
$ws.Rows.Item(4).Resize(2).Insert()

# Populate the two freshly-inserted cells in column A with the new
# Japanese "base hit-motion" strings.
$ws.Range("A4").Value = "ーーーーーー基礎ヤラレモーションーーーーーー"
$ws.Range("A5").Value = "0ちーんｗ 1倒れ 2ばたんきゅ 3溶ける 4燃える`n5射精 6吹っ飛ばす 7逃げこけ　8潰れる"

# Duplicate column A into column B for every row except the two newly
# inserted rows (4 and 5), which stay blank in column B.
$ws.Range("A1:A3").Copy()
$ws.Range("B1").PasteSpecial()

$ws.Range("A6:A53").Copy()
$ws.Range("B6").PasteSpecial()
